# CostValue, plus its digits
#
# - Recipe Qty (F10) drops from 50 to 20, which ripples into the Unit
#   Cost/"Jadi brp" formula in H10 (G10/F10).
# - The MIN() roll-up in H12 gets pointed at a bogus external/misspelled
#   sheet reference (h[1]Sheet!G8711), so it now evaluates to #VALUE!.
# - The view's selection/scroll position moves from J11 to J10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recipe Qty: 50 -> 20 (H10 = G10/F10 recalculates to 5 automatically)
$ws.Range("F10").Value = 20

# MIN(H10:H11) -> MIN(H10:h[1]Sheet!G8711), which errors out to #VALUE!
$ws.Range("H12").Formula = "=MIN(H10:h[1]Sheet!G8711)"

# Move the view's selection to J10 (was J11); also nudge the scrolled
# top-left cell toward C1 to match.
try {
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
$ws.Range("J10").Select()
